# Fix delimiter-detection bug: the "Measurements_*" field names for
# profundidade (depth), espessura (thickness), diametro (diameter) and
# peso (weight) were stored in lowercase, unlike their siblings
# (Measurements_Altura, Measurements_Largura) which are capitalized.
# Normalize the casing so downstream parsing that splits on the
# "Measurements_" delimiter/prefix behaves consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B74").Value2 = "Measurements_Profundidade"
$ws.Range("B75").Value2 = "Measurements_Espessura"
$ws.Range("B76").Value2 = "Measurements_Diametro"
$ws.Range("B77").Value2 = "Measurements_Peso"
